$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted at row 24 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 24-43 down to rows 25-44. The new row 24 repeats
# the same market/quality/price data as the (former) row 24 but is dated
# for the new week.
$ws.Rows("24:24").Insert()

$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44777
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112043
$ws.Range("G24").Value = "Pepino dulce"
$ws.Range("H24").Value = "Cultivar IV Región"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("N24").Value = '$/bandeja 18 kilos'
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 833
$ws.Range("Q24").Value = 18
$ws.Range("R24").Value = "Hortaliza"
